$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued numeric-looking Price strings are not auto-converted to numbers
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.140.54"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.798.59"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "337.03"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.4695"
$ws.Range("E7").Value = "  +24.54%  "
$ws.Range("D8").Value = "0.3707"
$ws.Range("E8").Value = "  +10.81%  "
$ws.Range("D9").Value = "45.25"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "0.07673"
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("D12").Value = "22.62"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "1.003"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "6.367"
$ws.Range("E14").Value = "  +3.81%  "
$ws.Range("D15").Value = "7.395"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").Value = "1.797.52"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "0.00001095"
$ws.Range("E17").Value = "  +3.90%  "
$ws.Range("D18").Value = "0.06754"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").Value = "82.54"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "17.42"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").Value = "6.421"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("D23").Value = "28.143.67"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "11.91"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("D25").Value = "2.412"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "20.85"
$ws.Range("E26").Value = "  +5.22%  "
$ws.Range("D27").Value = "2.389"
$ws.Range("E27").Value = "  +3.01%  "
$ws.Range("D28").Value = "151.65"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").Value = "2.003.94"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").Value = "133.67"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").Value = "1.261"
$ws.Range("D32").Value = "4.047"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").Value = "0.09698"
$ws.Range("E33").Value = "  +10.94%  "
$ws.Range("D34").Value = "5.931"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("D35").Value = "0.02384"
$ws.Range("E35").Value = "  +2.31%  "
$ws.Range("D36").Value = "12.19"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "0.2218"
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("D38").Value = "0.06352"
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("D39").Value = "0.6711"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").Value = "5.258"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").Value = "1.502"
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("D43").Value = "8.091"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "14.13"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "0.6159"
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("D47").Value = "3.853"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "130.43"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").Value = "2.063"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").Value = "1.182"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "0.07132"
$ws.Range("E51").Value = "  -0.52%  "

# Restore original (default) style on the price column
$priceRange.Style = "Normal"
